$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.190.24"
$ws.Range("E2").Value = "  -3.55%  "
$ws.Range("D3").Value = "1.750.85"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5818"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2709"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06605"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07513"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "1.750.89"
$ws.Range("E12").Value = "  -4.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.701"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6026"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "1.986.68"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "73.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008589"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -10.71%  "
$ws.Range("D18").Value = "28.159.96"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.288"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "204.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.654"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.995"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1232"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.403"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06052"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.395"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.715"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.699"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.656"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("E35").Value = "  -5.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6297"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.351"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.639"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.277"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01659"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("D41").Value = "1.129.50"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8605"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.009"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "1.898.25"
$ws.Range("E45").Value = "  -4.24%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000109"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.561"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.255"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.54%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05398"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4457"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.09%  "
